# Update F6:F13 progress values from 0.9 to 1 (Hold out and crossvalidation topic complete)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6:F13").Value = 1

# Match the selection change recorded in the saved file (active cell F6, selection F6:F13)
$ws.Range("F6:F13").Select()
